$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '42.217.10'
$ws.Range('E2').Value = '  +1.25%  '

# Row 3
$ws.Range('D3').Value = '2.172.47'
$ws.Range('E3').Value = '  -0.01%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '252.94'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +6.14%  '

# Row 6
$ws.Range('E6').Value = '  -1.09%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '73.38'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +1.04%  '

# Row 8
$ws.Range('E8').Value = '  -0.03%  '

# Row 9
$ws.Range('E9').Value = '  +0.10%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.12'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -0.08%  '

# Row 11
$ws.Range('E11').Value = '  -0.12%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '6.77'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +0.62%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.101'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.32%  '

# Row 14
$ws.Range('D14').Value = '2.499.16'
$ws.Range('E14').Value = '  +0.05%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.20'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.81%  '

# Row 16
$ws.Range('D16').Value = '2.197.11'
$ws.Range('E16').Value = '  +1.27%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.762'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.19%  '

# Row 18
$ws.Range('D18').Value = '42.106.36'
$ws.Range('E18').Value = '  +1.28%  '

# Row 19
$ws.Range('E19').Value = '  -0.43%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '70.58'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +0.77%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.84'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.81%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '226.36'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +0.04%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.52'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.56%  '

# Row 24
$ws.Range('E24').Value = '  +5.67%  '

# Row 25
$ws.Range('E25').Value = '  -0.21%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.44'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.76%  '

# Row 27
$ws.Range('E27').Value = '  +1.61%  '

# Row 28
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.19'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.05%  '

# Row 29
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.14'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.52%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.69'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +13.32%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '168.66'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -1.44%  '

# Row 32
$ws.Range('E32').Value = '  +0.70%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0804'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.82%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.11'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.90%  '

# Row 35
$ws.Range('E35').Value = '  -0.38%  '

# Row 36
$ws.Range('E36').Value = '  +4.00%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.22'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -2.13%  '

# Row 38
$ws.Range('E38').Value = '  +5.64%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '11.90'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.01%  '

# Row 40
$ws.Range('E40').Value = '  -2.62%  '

# Row 41
$ws.Range('E41').Value = '  +3.52%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '58.88'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -0.13%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.10'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -5.21%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '102.23'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +4.79%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.27'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.07%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.463'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +12.99%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0968'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.02%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.41'
$ws.Range('D48').Style = "Normal"

# Row 49
$ws.Range('E49').Value = '  +0.49%  '

# Row 50
$ws.Range('E50').Value = '  +0.54%  '

# Row 51
$ws.Range('E51').Value = '  +0.82%  '
